$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.809.75"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").Value = "3.493.98"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'581.17"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("D6").Value = "'161.82"
$ws.Range("E6").Value = "  +2.21%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'0.606"
$ws.Range("E8").Value = "  +9.39%  "

$ws.Range("D9").Value = "3.494.87"
$ws.Range("E9").Value = "  +0.99%  "

$ws.Range("D10").Value = "'7.33"
$ws.Range("E10").Value = "  -2.95%  "

$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("D13").Value = "4.095.12"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("E14").Value = "  -0.89%  "

$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").Value = "'28.89"
$ws.Range("E16").Value = "  +4.17%  "

$ws.Range("D17").Value = "65.791.85"
$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("D18").Value = "3.520.85"
$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("D19").Value = "'6.47"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").Value = "'14.34"
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("D21").Value = "'392.29"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").Value = "'8.30"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("D24").Value = "'73.70"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("D27").Value = "'9.76"
$ws.Range("E27").Value = "  +0.64%  "

$ws.Range("D28").Value = "'0.179"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").Value = "'6.40"
$ws.Range("E30").Value = "  +3.32%  "

$ws.Range("D31").Value = "'1.46"
$ws.Range("E31").Value = "  +5.61%  "

$ws.Range("E32").Value = "  +1.40%  "

$ws.Range("D33").Value = "'23.80"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "'6.56"
$ws.Range("E34").Value = "  -1.36%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'7.20"
$ws.Range("E36").Value = "  +1.87%  "

$ws.Range("E37").Value = "  +5.36%  "

$ws.Range("D38").Value = "'162.83"
$ws.Range("E38").Value = "  +1.23%  "

$ws.Range("E39").Value = "  +5.17%  "

$ws.Range("D40").Value = "3.109.20"
$ws.Range("E40").Value = "  +6.10%  "

$ws.Range("D41").Value = "'0.0775"
$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("D42").Value = "'27.38"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").Value = "'0.0324"
$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("E44").Value = "  +2.12%  "

$ws.Range("D45").Value = "'43.15"
$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").Value = "'25.99"
$ws.Range("E47").Value = "  +7.97%  "

$ws.Range("E48").Value = "  +2.97%  "

$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'314.24"
$ws.Range("E50").Value = "  +6.32%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.74"
$ws.Range("E51").Value = "  +2.72%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
